$d = $word.ActiveDocument

# 1. Insert "Databricks | " between "AWS | " and "Snowflake " in the skills list.
$d.Content.Find.Execute("AWS | Snowflake", $true, $false, $false, $false, $false, $true, 1, $false, "AWS | Databricks | Snowflake", 2)

# 2. Remove " Plotly" so that "Qlik | Plotly | Dbeaver" becomes "Qlik | Dbeaver"
#    (the run/proofErr pair that used to spell-check "Plotly" disappears, and the
#    following "Dbeaver" run - with its own proofErr wrapper - is untouched).
$d.Content.Find.Execute("Qlik | Plotly |", $true, $false, $false, $false, $false, $true, 1, $false, "Qlik |", 2)

# 3. Shorten "MS Office Suite" to "MS Office".
$d.Content.Find.Execute("MS Office Suite", $true, $false, $false, $false, $false, $true, 1, $false, "MS Office", 2)

# 4. Normalize the two adjacent runs describing the AI roadmap bullet into one
#    (self-replace forces Word to merge the runs into a single run, matching
#    the canonical output, without altering the visible text).
$d.Content.Find.Execute("Developed a 12-month AI roadmap prioritizing high-value use cases across Operations, Research, and Manufacturing.", $true, $false, $false, $false, $false, $true, 1, $false, "Developed a 12-month AI roadmap prioritizing high-value use cases across Operations, Research, and Manufacturing.", 2)
